$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F5").Value = 1196
$ws.Range("F6").Value = 9363
$ws.Range("F10").Value = 324
$ws.Range("F11").Value = 5719
$ws.Range("F13").Value = 15
$ws.Range("F14").Value = 6585
$ws.Range("F16").Value = 459
$ws.Range("F17").Value = 437
$ws.Range("F18").Value = 640
$ws.Range("F19").Value = 339
$ws.Range("F21").Value = 219
$ws.Range("F24").Value = 110
$ws.Range("F25").Value = 10648
$ws.Range("F26").Value = 92
$ws.Range("F28").Value = 2031
$ws.Range("F29").Value = 2481
$ws.Range("F30").Value = 47
$ws.Range("F32").Value = 2330
$ws.Range("F37").Value = 317
$ws.Range("F38").Value = 1467
$ws.Range("F40").Value = 5443
$ws.Range("F41").Value = 1222
$ws.Range("F42").Value = 745
$ws.Range("F43").Value = 134
$ws.Range("F46").Value = 1088
$ws.Range("F47").Value = 1426
$ws.Range("F49").Value = 1110

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F11").Value = 199
$ws.Range("F16").Value = 101
$ws.Range("F19").Value = 915
$ws.Range("F20").Value = 17

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 70
$ws.Range("F3").Value = 144

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 1196
$ws.Range("F6").Value = 9363
$ws.Range("F9").Value = 324
$ws.Range("F10").Value = 144
$ws.Range("F13").Value = 5719
$ws.Range("F14").Value = 5719
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 6585
$ws.Range("F18").Value = 6585
$ws.Range("F20").Value = 459
$ws.Range("F21").Value = 437
$ws.Range("F22").Value = 640
$ws.Range("F23").Value = 339
$ws.Range("F25").Value = 219
$ws.Range("F28").Value = 199
$ws.Range("F29").Value = 10648
$ws.Range("F30").Value = 92
$ws.Range("F32").Value = 2031
$ws.Range("F33").Value = 2481
$ws.Range("F34").Value = 2330
$ws.Range("F37").Value = 317
$ws.Range("F38").Value = 1467
$ws.Range("F40").Value = 5443
$ws.Range("F41").Value = 17
$ws.Range("F42").Value = 1222
$ws.Range("F43").Value = 745
$ws.Range("F44").Value = 134
$ws.Range("F47").Value = 1088
$ws.Range("F49").Value = 1426
$ws.Range("F51").Value = 1110
